# Apply weekly price/date updates to the "Ajo" worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Final values per row for columns D (Fecha), J (Volumen), K (Precio minimo),
# L (Precio maximo), M (Precio promedio ponderado), P (Precio $/Kg)
$rows = @{
    2  = @{ D = 44441; J = 300; K = 15000; L = 16000; M = 15500; P = 1550 }
    3  = @{ D = 44524; J = 200; K = 20000; L = 21000; M = 20500; P = 2050 }
    4  = @{ D = 44263; J = 300; K = 15000; L = 16000; M = 15500; P = 1550 }
    5  = @{ D = 44265; J = 200; K = 15000; L = 16000; M = 15500; P = 1550 }
    6  = @{ D = 44204; J = 400; K = 10000; L = 11000; M = 10500; P = 1050 }
    7  = @{ D = 44330; J = 300; K = 13000; L = 14000; M = 13500; P = 1350 }
    8  = @{ D = 44377; J = 650; K = 14000; L = 15000; M = 14538; P = 1454 }
    9  = @{ D = 44460; J = 300; K = 15000; L = 16000; M = 15500; P = 1550 }
    10 = @{ D = 44358; J = 300; K = 14000; L = 15000; M = 14500; P = 1450 }
    11 = @{ D = 44291; J = 200; K = 13000; L = 14000; M = 13500; P = 1350 }
    12 = @{ D = 44218; J = 320; K = 10000; L = 11000; M = 10500; P = 1050 }
    13 = @{ D = 44428; J = 300; K = 15000; L = 16000; M = 15500; P = 1550 }
    14 = @{ D = 44406; J = 400; K = 14000; L = 15000; M = 14500; P = 1450 }
    15 = @{ D = 44580; J = 200; K = 18000; L = 20000; M = 19000; P = 1900 }
    16 = @{ D = 44644; J = 300; K = 20000; L = 21000; M = 20500; P = 2050 }
    17 = @{ D = 44547; J = 300; K = 19000; L = 20000; M = 19500; P = 1950 }
    18 = @{ D = 44160; J = 360; K = 10000; L = 11000; M = 10500; P = 1050 }
}

foreach ($r in $rows.Keys) {
    $vals = $rows[$r]
    $ws.Cells.Item($r, 4).Value  = $vals.D   # D - Fecha
    $ws.Cells.Item($r, 10).Value = $vals.J   # J - Volumen
    $ws.Cells.Item($r, 11).Value = $vals.K   # K - Precio minimo
    $ws.Cells.Item($r, 12).Value = $vals.L   # L - Precio maximo
    $ws.Cells.Item($r, 13).Value = $vals.M   # M - Precio promedio ponderado
    $ws.Cells.Item($r, 16).Value = $vals.P   # P - Precio $/Kg
}
